$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "Power Regulator Circuit"
$ws.Range("G5").Value = "Toggle Switch"

$ws.Range("G6").Select()
